$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$players = @(
    @("Usman Khawaja", 1, "Batsman", "", "Australia"),
    @("Alex Carey", 4, "Batsman", "Wicket-Keeper", "Australia"),
    @("Aaron Finch", 5, "Batsman", "Captain", "Australia"),
    @("Tim Paine", 7, "Batsman", "Captain", "Australia"),
    @("Shaun Marsh", 9, "Batsman", "", "Australia"),
    @("Matthew Wade", 13, "Batsman", "Wicket-Keeper", "Australia"),
    @("Marcus Harris", 14, "Batsman", "Wicket-Keeper", "Australia"),
    @("Marcus Stoinis", 17, "All Rounder", "", "Australia"),
    @("D'Arcy Short", 23, "All Rounder", "", "Australia"),
    @("Pat Cummins", 30, "Bowler", "", "Australia"),
    @("David Warner", 31, "Batsman", "Vice-Captain", "Australia"),
    @("Glenn Maxwell", 32, "All Rounder", "", "Australia"),
    @("Marnus Labuschagne", 33, "Batsman", "", "Australia"),
    @("Ashton Agar", 46, "All Rounder", "", "Australia"),
    @("Steve Smith", 49, "Batsman", "Captain", "Australia"),
    @("Mitchell Stark", 56, "Bowler", "", "Australia"),
    @("Jhye Richardson", 60, "Bowler", "", "Australia")
)

$startRow = 17

# Write column E (team) first for all rows so "Australia" is registered
# as a shared string before the player names, matching the author's order.
for ($i = 0; $i -lt $players.Count; $i++) {
    $row = $startRow + $i
    $p = $players[$i]
    $ws.Cells.Item($row, 5).Value = $p[4]
}

for ($i = 0; $i -lt $players.Count; $i++) {
    $row = $startRow + $i
    $p = $players[$i]
    $ws.Cells.Item($row, 1).Value = $p[0]
    $ws.Cells.Item($row, 2).Value = $p[1]
    $ws.Cells.Item($row, 3).Value = $p[2]
    if ($p[3] -ne "") {
        $ws.Cells.Item($row, 4).Value = $p[3]
    }
}

$ws.Range("B3").Select()
